# Initial Data File Update
# Adds six new transaction rows (209-214) to the "Transacciones" sheet,
# continuing the running-balance ledger that ends at row 208.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Transacciones")

# ---------------------------------------------------------------------
# Row 209 - Panquecitos Gotas (snack), paid with Tarjeta Santander
# ---------------------------------------------------------------------
$ws.Range("A209").Value = 43604
$ws.Range("B209").Value = 19.5
$ws.Range("C209").Value = "Panquecitos Gotas"
$ws.Range("D209").Value = "Golosinas"
$ws.Range("E209").Value = "Gasto"
$ws.Range("F209").Value = "Tarjeta Santander"
$ws.Range("G209").Value = "Extra"
$ws.Range("K209").Value = 4237.18
$ws.Range("L209").Formula = "=L208-B209"
$ws.Range("M209").Value = 127
$ws.Range("N209").Formula = "=SUM(K209:M209)"
$ws.Range("O209").Formula = "=N209-4000"
$ws.Range("P209").Formula = "=O209-Ahorros!`$E`$4"

# ---------------------------------------------------------------------
# Row 210 - Lavanderia, paid with Tarjeta Banamex
# ---------------------------------------------------------------------
$ws.Range("A210").Value = 43604
$ws.Range("B210").Value = 65
$ws.Range("C210").Value = "Lavandería"
$ws.Range("D210").Value = "Lavandería"
$ws.Range("E210").Value = "Gasto"
$ws.Range("F210").Value = "Tarjeta Banamex"
$ws.Range("G210").Value = "Lavandería"
$ws.Range("K210").Formula = "=K209-B210"
$ws.Range("L210").Value = 5458.26
$ws.Range("M210").Value = 127
$ws.Range("N210").Formula = "=SUM(K210:M210)"
$ws.Range("O210").Formula = "=N210-4000"
$ws.Range("P210").Formula = "=O210-Ahorros!`$E`$4"

# ---------------------------------------------------------------------
# Row 211 - Alitas Hot Wings, paid with Tarjeta Santander
# ---------------------------------------------------------------------
$ws.Range("A211").Value = 43605
$ws.Range("B211").Value = 181
$ws.Range("C211").Value = "Alitas Hot Wings"
$ws.Range("D211").Value = "Comida"
$ws.Range("E211").Value = "Gasto"
$ws.Range("F211").Value = "Tarjeta Santander"
$ws.Range("G211").Value = "Hot Wings"
$ws.Range("K211").Value = 4172.18
$ws.Range("L211").Formula = "=L210-B211"
$ws.Range("M211").Value = 127
$ws.Range("N211").Formula = "=SUM(K211:M211)"
$ws.Range("O211").Formula = "=N211-4000"
$ws.Range("P211").Formula = "=O211-Ahorros!`$E`$4"

# ---------------------------------------------------------------------
# Row 212 - Propina Hot Wings, paid in cash (Efectivo)
# ---------------------------------------------------------------------
$ws.Range("A212").Value = 43605
$ws.Range("B212").Value = 19
$ws.Range("C212").Value = "Propina Hot Wings"
$ws.Range("D212").Value = "Propina"
$ws.Range("E212").Value = "Gasto"
$ws.Range("F212").Value = "Efectivo"
$ws.Range("G212").Value = "Hot Wings"
$ws.Range("K212").Value = 4172.18
$ws.Range("L212").Value = 5277.26
$ws.Range("M212").Formula = "=M211-B212"
$ws.Range("N212").Formula = "=SUM(K212:M212)"
$ws.Range("O212").Formula = "=N212-4000"
$ws.Range("P212").Formula = "=O212-Ahorros!`$E`$4"

# ---------------------------------------------------------------------
# Row 213 - Uber, paid in cash (Efectivo)
# ---------------------------------------------------------------------
$ws.Range("A213").Value = 43605
$ws.Range("B213").Value = 10
$ws.Range("C213").Value = "Uber"
$ws.Range("D213").Value = "Transporte"
$ws.Range("E213").Value = "Gasto"
$ws.Range("F213").Value = "Efectivo"
$ws.Range("G213").Value = "NA"
$ws.Range("K213").Value = 4172.18
$ws.Range("L213").Value = 5277.26
$ws.Range("M213").Formula = "=M212-B213"
$ws.Range("N213").Formula = "=SUM(K213:M213)"
$ws.Range("O213").Formula = "=N213-4000"
$ws.Range("P213").Formula = "=O213-Ahorros!`$E`$4"

# ---------------------------------------------------------------------
# Row 214 - Ahorro en Alcancia, paid in cash (Efectivo)
# ---------------------------------------------------------------------
$ws.Range("A214").Value = 43605
$ws.Range("B214").Value = 10
$ws.Range("C214").Value = "Ahorro en Alcancía"
$ws.Range("D214").Value = "Ahorro"
$ws.Range("E214").Value = "Gasto"
$ws.Range("F214").Value = "Efectivo"
$ws.Range("G214").Value = "Alcancía"
$ws.Range("K214").Value = 4172.18
$ws.Range("L214").Value = 5277.26
$ws.Range("M214").Formula = "=M213-B214"
$ws.Range("N214").Formula = "=SUM(K214:M214)"
$ws.Range("O214").Formula = "=N214-4000"
$ws.Range("P214").Formula = "=O214-Ahorros!`$E`$4"

# ---------------------------------------------------------------------
# Formatting: carry the date format (column A) and the highlighted
# "ahorro" style (column P) down from the last existing row (208) onto
# the new rows, the same way dragging the fill handle would.
# ---------------------------------------------------------------------
$ws.Range("A208").Copy()
$ws.Range("A209:A214").PasteSpecial(-4122)

$ws.Range("P208").Copy()
$ws.Range("P209:P214").PasteSpecial(-4122)

$excel.CutCopyMode = 0

# ---------------------------------------------------------------------
# Leave the selection where the user's edit session ended.
# ---------------------------------------------------------------------
$ws.Range("Q214").Select()
